$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 8 (eh_direita=TRUE, favoravel, "Racismo reverso não existe", 1).
# This shifts rows 9-15 up to become rows 8-14 and updates the sheet dimension
# from A1:D15 to A1:D14 automatically.
$ws.Rows.Item(8).Delete()

# Update the follower counts ("n" column) that changed, using the row numbers
# as they stand after the deletion/shift above.
$ws.Range("D3").Value = 6
$ws.Range("D4").Value = 44
$ws.Range("D5").Value = 8
$ws.Range("D6").Value = 6

$ws.Range("D9").Value = 6
$ws.Range("D10").Value = 27
$ws.Range("D11").Value = 3
$ws.Range("D12").Value = 10
$ws.Range("D14").Value = 171
